# Generate Report for Handoff
# Adds two new files (6fd19b7e-... and 7c2f463d-...) to the localization
# status report: one new row per file on the "Overview" sheet and one new
# row per file on each of the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3 "Overview")
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)

$tblOverview.ListRows.Add()
$wsOverview.Range("A6").Value = "6fd19b7e-6d10-465c-aeb4-5d726acf5479.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = "'"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-09-01 14:56:57"
$wsOverview.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/067b0f6ebdfde4beda5ff522e6bd954bf007ec1e/e2e/6fd19b7e-6d10-465c-aeb4-5d726acf5479.md", "", "", "e2e\6fd19b7e-6d10-465c-aeb4-5d726acf5479.md")

$tblOverview.ListRows.Add()
$wsOverview.Range("A7").Value = "7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = "'"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-09-01 14:56:57"
$wsOverview.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/067b0f6ebdfde4beda5ff522e6bd954bf007ec1e/e2e/7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md", "", "", "e2e\7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1 "zh-cn")
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item(1)

$tblZhCn.ListRows.Add()
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "'ht"
$wsZhCn.Range("F6").Value = "'False"
$wsZhCn.Range("G6").Value = "6fd19b7e-6d10-465c-aeb4-5d726acf5479.480464cd7ab3d79e0af08f4c1d129a2bb115a1c9.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-09-01 14:56:51"
$wsZhCn.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I6").Value = "'"
$wsZhCn.Range("J6").Value = "'"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L6").Value = "'"
$wsZhCn.Range("M6").Value = "'True"
$wsZhCn.Range("N6").Value = "'"
$wsZhCn.Range("O6").Value = "'False"
$wsZhCn.Range("P6").Value = "'"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/067b0f6ebdfde4beda5ff522e6bd954bf007ec1e/e2e/6fd19b7e-6d10-465c-aeb4-5d726acf5479.md", "", "", "6fd19b7e-6d10-465c-aeb4-5d726acf5479.md")

$tblZhCn.ListRows.Add()
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "'ht"
$wsZhCn.Range("F7").Value = "'False"
$wsZhCn.Range("G7").Value = "7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.9d45f928b8b775dd97ee846c696742e71612ecdd.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-09-01 14:56:51"
$wsZhCn.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I7").Value = "'"
$wsZhCn.Range("J7").Value = "'"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L7").Value = "'"
$wsZhCn.Range("M7").Value = "'True"
$wsZhCn.Range("N7").Value = "'"
$wsZhCn.Range("O7").Value = "'False"
$wsZhCn.Range("P7").Value = "'"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/067b0f6ebdfde4beda5ff522e6bd954bf007ec1e/e2e/7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md", "", "", "7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md")

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2 "de-de")
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item(1)

$tblDeDe.ListRows.Add()
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "'ht"
$wsDeDe.Range("F6").Value = "'False"
$wsDeDe.Range("G6").Value = "6fd19b7e-6d10-465c-aeb4-5d726acf5479.480464cd7ab3d79e0af08f4c1d129a2bb115a1c9.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-09-01 14:56:57"
$wsDeDe.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I6").Value = "'"
$wsDeDe.Range("J6").Value = "'"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L6").Value = "'"
$wsDeDe.Range("M6").Value = "'True"
$wsDeDe.Range("N6").Value = "'"
$wsDeDe.Range("O6").Value = "'False"
$wsDeDe.Range("P6").Value = "'"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/067b0f6ebdfde4beda5ff522e6bd954bf007ec1e/e2e/6fd19b7e-6d10-465c-aeb4-5d726acf5479.md", "", "", "6fd19b7e-6d10-465c-aeb4-5d726acf5479.md")

$tblDeDe.ListRows.Add()
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "'ht"
$wsDeDe.Range("F7").Value = "'False"
$wsDeDe.Range("G7").Value = "7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.9d45f928b8b775dd97ee846c696742e71612ecdd.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-09-01 14:56:57"
$wsDeDe.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I7").Value = "'"
$wsDeDe.Range("J7").Value = "'"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L7").Value = "'"
$wsDeDe.Range("M7").Value = "'True"
$wsDeDe.Range("N7").Value = "'"
$wsDeDe.Range("O7").Value = "'False"
$wsDeDe.Range("P7").Value = "'"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/067b0f6ebdfde4beda5ff522e6bd954bf007ec1e/e2e/7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md", "", "", "7c2f463d-c0cb-4286-8956-ba0fc5dd76f5.md")
